$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Покупка билета) ---
$ws.Range("C2").Value = 25
$ws.Range("D2").Formula = "=60/C2"
$ws.Range("H2").Formula = "=D2*E2*G2"

# --- Row 3 (Логин и логаут) ---
$ws.Range("D3").Formula = "=60/C3"
$ws.Range("H3").Formula = "=D3*E3*G3"

# --- Row 4 (Удаление брони) ---
# C4's cell format carries an unusual "quote prefix" style flag that a
# plain Value write would drop, so stash the current formatting, update
# the value/formula, then restore the original formatting.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C4").Value = 19

$ws.Range("K1").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K1").Clear() | Out-Null

$ws.Range("D4").Formula = "=60/C4"
$ws.Range("H4").Formula = "=D4*E4*G4"

# --- Row 5 (Поиск билета) ---
$ws.Range("C5").Value = 32
$ws.Range("D5").Formula = "=60/C5"
$ws.Range("H5").Formula = "=D5*E5*G5"

# --- Row 6 (Просмотр текущих бронирований) ---
$ws.Range("C6").Value = 15
$ws.Range("D6").Formula = "=60/C6"
$ws.Range("H6").Formula = "=D6*E6*G6"

# --- Header text update: intensity is now computed for ALL Vus, not one ---
$ws.Range("H1").Value = "Расчётная интенсивность операций за 1 час всеми Vus"

# --- Selected cell moved from H10 to D8 ---
$ws.Range("D8").Select() | Out-Null
